$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Coach G" row (id 7) to the dataset
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Coach G"
$ws.Range("C8").Value = "G"

# Update the current selection to match the edited workbook
[void]$ws.Range("J10").Select()
